# Add test data for pduration. Begins #62
#
# Adds a new "pduration" worksheet (after the existing "rri" sheet) with a
# Table2 listobject (rate/pv/fv/pduration columns) and PDURATION() test
# rows, mirroring the existing "rri"/RRI() sheet/table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# New worksheet, placed after the "rri" sheet.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "pduration"

# Header row.
$ws2.Range("A1").Value = "rate"
$ws2.Range("B1").Value = "pv"
$ws2.Range("C1").Value = "fv"
$ws2.Range("D1").Value = "pduration"

# Row 2 - standalone formula (not part of the shared group below).
$ws2.Cells.Item(2, 1).Value = 0.015309470499731193
$ws2.Cells.Item(2, 2).Value = -5
$ws2.Cells.Item(2, 3).Value = -6
$ws2.Range("D2").Formula = "=_xlfn.PDURATION(A2,B2,C2)"

# Rows 3-24 input data (rate, pv, fv).
$ws2.Cells.Item(3, 1).Value = -1
$ws2.Cells.Item(3, 2).Value = -5
$ws2.Cells.Item(3, 3).Value = 0

$ws2.Cells.Item(4, 1).Value = 0
$ws2.Cells.Item(4, 2).Value = -1
$ws2.Cells.Item(4, 3).Value = -1

$ws2.Cells.Item(5, 1).Value = 0
$ws2.Cells.Item(5, 2).Value = 300
$ws2.Cells.Item(5, 3).Value = 300

$ws2.Cells.Item(6, 1).Value = 0.1
$ws2.Cells.Item(6, 2).Value = 0
$ws2.Cells.Item(6, 3).Value = 100

$ws2.Cells.Item(7, 1).Value = 0.1
$ws2.Cells.Item(7, 2).Value = 100
$ws2.Cells.Item(7, 3).Value = 0

$ws2.Cells.Item(8, 1).Value = 0.02426318074098921
$ws2.Cells.Item(8, 2).Value = 300
$ws2.Cells.Item(8, 3).Value = 400

$ws2.Cells.Item(9, 1).Value = 0.24092317318260137
$ws2.Cells.Item(9, 2).Value = 300
$ws2.Cells.Item(9, 3).Value = 4000

$ws2.Cells.Item(10, 1).Value = 0.50341274654387536
$ws2.Cells.Item(10, 2).Value = 300
$ws2.Cells.Item(10, 3).Value = 40000

$ws2.Cells.Item(11, 1).Value = 0.012058882052318642
$ws2.Cells.Item(11, 2).Value = 300
$ws2.Cells.Item(11, 3).Value = 400

$ws2.Cells.Item(12, 1).Value = 0.11396731243901459
$ws2.Cells.Item(12, 2).Value = 300
$ws2.Cells.Item(12, 3).Value = 4000

$ws2.Cells.Item(13, 1).Value = 0.22613732776711237
$ws2.Cells.Item(13, 2).Value = 300
$ws2.Cells.Item(13, 3).Value = 40000

$ws2.Cells.Item(14, 1).Value = 0.0075993101546305564
$ws2.Cells.Item(14, 2).Value = 300
$ws2.Cells.Item(14, 3).Value = 400

$ws2.Cells.Item(15, 1).Value = 0.070541853470322824
$ws2.Cells.Item(15, 2).Value = 300
$ws2.Cells.Item(15, 3).Value = 4000

$ws2.Cells.Item(16, 1).Value = 0.13741628093790048
$ws2.Cells.Item(16, 2).Value = 300
$ws2.Cells.Item(16, 3).Value = 40000

$ws2.Cells.Item(17, 1).Value = 0.98822504304098735
$ws2.Cells.Item(17, 2).Value = 10000
$ws2.Cells.Item(17, 3).Value = 2441880

$ws2.Cells.Item(18, 1).Value = 0.046635139392105618
$ws2.Cells.Item(18, 2).Value = 5000
$ws2.Cells.Item(18, 3).Value = 6000

$ws2.Cells.Item(19, 1).Value = 0.18920711500272103
$ws2.Cells.Item(19, 2).Value = 5000
$ws2.Cells.Item(19, 3).Value = 10000

$ws2.Cells.Item(20, 1).Value = 0.10000000000000009
$ws2.Cells.Item(20, 2).Value = 250
$ws2.Cells.Item(20, 3).Value = 275

$ws2.Cells.Item(21, 1).Value = 0.41421356237309492
$ws2.Cells.Item(21, 2).Value = 250
$ws2.Cells.Item(21, 3).Value = 500

$ws2.Cells.Item(22, 1).Value = 0.5211809843045565
$ws2.Cells.Item(22, 2).Value = 250
$ws2.Cells.Item(22, 3).Value = 880

$ws2.Cells.Item(23, 1).Value = 0.025000000000000001
$ws2.Cells.Item(23, 2).Value = 2000
$ws2.Cells.Item(23, 3).Value = 2200

# Row 24 - rate is itself a formula.
$ws2.Range("A24").Formula = "=0.025/12"
$ws2.Cells.Item(24, 2).Value = 1000
$ws2.Cells.Item(24, 3).Value = 1200

# Rows 3-24 share one formula definition (fill down from D3).
$ws2.Range("D3:D24").Formula = "=_xlfn.PDURATION(A3,B3,C3)"

# ---------------------------------------------------------------------
# Turn A1:D24 into a table, matching Table1 on the "rri" sheet.
# ---------------------------------------------------------------------
$lo = $ws2.ListObjects.Add(1, $ws2.Range("A1:D24"), $null, 1)
$lo.Name = "Table2"

# Column width tweak on column D (seen in the diff for this sheet).
$ws2.Columns.Item(4).ColumnWidth = 11.08984375

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the diff.
# ---------------------------------------------------------------------
$ws1.Range("B14:D32").Select()
$ws1.Activate()

$ws2.Range("A24").Select()
$ws2.Activate()
